$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for columns B,C,E,F,G,I,J,M,N across rows 2-25
# (columns A,D,H,K,L,O are unchanged and left as-is)
$columns = @(2, 3, 5, 6, 7, 9, 10, 13, 14)  # B, C, E, F, G, I, J, M, N

$data = @{
    2 = @(16.2005927474062, 12.83023013918517, 16.07902887608318, 43.5165815734676, 3.663456794202379, 24.9964988682794, 8.665293179902653, 19.80173307187955, 18.6529744591511)
    3 = @(15.6759274972307, 12.28764008112553, 16.08517181177352, 43.32905868479114, 3.667355850194007, 24.99908159120005, 8.693256506517024, 19.66557319454697, 18.7266668312664)
    4 = @(15.34962376500426, 11.94617556115307, 16.09038227704462, 43.22584265150275, 3.669872218079639, 25.00773384570566, 8.711439869489308, 19.5860971666723, 18.77390899403686)
    5 = @(15.21583849444901, 11.80517162510868, 16.09286860904798, 43.18680466177351, 3.670928540761172, 25.01302916232764, 8.719104960014782, 19.5547754989123, 18.79366354076593)
    6 = @(15.19358147493964, 11.78165345450484, 16.09330342442296, 43.18050568265569, 3.67110581110643, 25.01401507543612, 8.720393166346442, 19.54963968643061, 18.79697418787498)
    7 = @(15.34782246772187, 11.94428110940004, 16.09041433703649, 43.22530390039532, 3.669886338817092, 25.00779810665521, 8.711542209808561, 19.58567040242956, 18.77417337244459)
    8 = @(16.02068483839479, 12.64500333961482, 16.08084912685078, 43.449465262028, 3.664775871782211, 24.99591859975867, 8.674724708161554, 19.75394585703321, 18.67797041778694)
    9 = @(17.29785528692536, 13.94415644567042, 16.07345148736337, 43.98235186956529, 3.655719383024926, 25.02897543723794, 8.610555797793303, 20.11537472140938, 18.50507929380897)
    10 = @(18.19914009365573, 14.84264268269057, 16.07486654948321, 44.42872013695941, 3.649646202862412, 25.08792932843413, 8.568287298162907, 20.39823440288206, 18.38757246192473)
    11 = @(18.59908401615674, 15.23750669517003, 16.0769809517887, 44.64318143809269, 3.647007776077434, 25.12231725572682, 8.550113341103181, 20.53025445074302, 18.33616215523557)
    12 = @(18.74894050459899, 15.38492136343399, 16.07799158252791, 44.72598348746543, 3.646026418196503, 25.13642859850132, 8.543382689476733, 20.58069126371143, 18.31698693335132)
    13 = @(18.71673928424295, 15.35326859931313, 16.07776461231094, 44.7080806900011, 3.646236983268951, 25.13334101793037, 8.544825523156424, 20.56980966381151, 18.32110365939988)
    14 = @(18.61144552619176, 15.24967752789509, 16.0770598981599, 44.64996204026057, 3.64692668396001, 25.12345637126203, 8.549556572294978, 20.53439524242983, 18.33457873793808)
    15 = @(18.54673841912383, 15.18594677569352, 16.07665553618056, 44.61456824199732, 3.647351454602257, 25.11754361278548, 8.552474192850672, 20.51275954492225, 18.34287070146789)
    16 = @(18.17278748708045, 14.81654792365173, 16.07475783254481, 44.41492991072919, 3.649821121420676, 25.08583442172866, 8.569496208019833, 20.38967072462428, 18.39097327671377)
    17 = @(17.94069507470924, 14.5862933495347, 16.07396925685075, 44.29534515300284, 3.651367933768582, 25.06832165815028, 8.580208548465103, 20.31499025893534, 18.42100534217198)
    18 = @(17.80626299122799, 14.45255814011229, 16.07365421124814, 44.22763997649941, 3.652269325320054, 25.05896129001626, 8.586469237868693, 20.27235324306887, 18.43847146662235)
    19 = @(17.76059014876601, 14.40705851613003, 16.0735713793897, 44.20490257867561, 3.65257653507154, 25.0559143549817, 8.588606048481283, 20.257972691903, 18.44441829079314)
    20 = @(17.96549997812941, 14.61093968406609, 16.07403887648107, 44.3079640600773, 3.651202062269138, 25.07011216886327, 8.579057932359406, 20.32290756052422, 18.41778846589377)
    21 = @(18.64241720191192, 15.28016290668491, 16.07726120417306, 44.66699015317507, 3.646723621150424, 25.12633016567087, 8.54816283987304, 20.54478557666548, 18.330612845156)
    22 = @(19.07547574672768, 15.70518697687081, 16.08059020159242, 44.91087799684522, 3.64390014918368, 25.16942108268419, 8.528853706343567, 20.69236402329953, 18.27534434015347)
    23 = @(18.84524407002672, 15.47950854415765, 16.07870206450232, 44.7798820195999, 3.645397662078333, 25.14584175018079, 8.539078645652623, 20.61337616014886, 18.30468649135007)
    24 = @(17.95428878177488, 14.5998012993748, 16.07400697050372, 44.30225579133597, 3.651277015049793, 25.0693004739385, 8.579577808112242, 20.31932721813216, 18.41924219170978)
    25 = @(16.95810742033207, 13.60188162670681, 16.07424417767032, 43.82839566963497, 3.658066881993338, 25.01396242516192, 8.627057526331578, 20.01441808097795, 18.55017311316885)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Cells.Item($row, $columns[$i]).Value = $values[$i]
    }
}
